# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# This script:
#  - Removes the two worker rows that were taken out of the statement
#    (LAURA ANDREA VASQUEZ CASTELLAR / period 1704, and HERMEN GOMEZ HERRERA / period 2503),
#    keeping only the JOSE BERTEL MELGAREJO row.
#  - Updates the summary totals (VALOR MORA, Cant. Trabajadores, Cant. Periodos)
#    to reflect the reduced data set.
#  - Adjusts column D's width to fit the remaining (shorter) worker name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the aggregate "VALOR MORA" total (was the sum for 3 workers, now just 1)
$ws.Range("E11").Value2 = 31816

# Update the worker / period counters
$ws.Range("C13").Value2 = 1
$ws.Range("F13").Value2 = 1

# Remove the rows for the two workers that are no longer part of this statement
$ws.Rows("17:19").Delete()

# Column D no longer needs to fit the long removed name; shrink to fit what's left
$ws.Columns("D").ColumnWidth = 24
